$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date value in C2 (45656 -> 45678)
$ws.Range("C2").Value = 45678

# Update the active selection to E2
$ws.Range("E2").Select()
